$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "40×76=3040"; New = "90×93=8370" },
    @{ Old = "65×85=5525"; New = "17×41=697" },
    @{ Old = "47×92=4324"; New = "61×34=2074" },
    @{ Old = "46×74=3404"; New = "99×43=4257" },
    @{ Old = "78×60=4680"; New = "91×68=6188" },
    @{ Old = "52×74=3848"; New = "39×57=2223" },
    @{ Old = "40×79=3160"; New = "89×31=2759" },
    @{ Old = "77×88=6776"; New = "68×73=4964" },
    @{ Old = "33×45=1485"; New = "29×77=2233" },
    @{ Old = "36×57=2052"; New = "92×54=4968" },
    @{ Old = "61×66=4026"; New = "18×33=594" },
    @{ Old = "27×77=2079"; New = "16×95=1520" },
    @{ Old = "70×36=2520"; New = "44×47=2068" },
    @{ Old = "22×50=1100"; New = "72×40=2880" },
    @{ Old = "96×90=8640"; New = "87×25=2175" },
    @{ Old = "81×53=4293"; New = "21×27=567" },
    @{ Old = "56×64=3584"; New = "46×12=552" },
    @{ Old = "85×28=2380"; New = "69×69=4761" },
    @{ Old = "40×11=440";  New = "88×79=6952" },
    @{ Old = "95×32=3040"; New = "30×80=2400" },
    @{ Old = "12×65=780";  New = "97×66=6402" },
    @{ Old = "40×31=1240"; New = "60×88=5280" },
    @{ Old = "61×42=2562"; New = "98×95=9310" },
    @{ Old = "86×41=3526"; New = "37×66=2442" },
    @{ Old = "54×99=5346"; New = "57×39=2223" }
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
